# Auto-generated Excel COM-interop edit script
# Applies the numeric corrections to the Leve profit-tracking sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 3600
$ws.Range("J64").Value = 3600
$ws.Range("L64").Value = 3600
$ws.Range("N64").Value = -4096
# Row 67
$ws.Range("H67").Value = 3600
$ws.Range("J67").Value = 3600
$ws.Range("L67").Value = 3600
$ws.Range("N67").Value = -5316
# Row 94
$ws.Range("H94").Value = 8602
$ws.Range("I94").Value = 8602
$ws.Range("K94").Value = 8602
$ws.Range("M94").Value = -8151
# Row 99
$ws.Range("H99").Value = 200
$ws.Range("I99").Value = 200
$ws.Range("K99").Value = 600
$ws.Range("M99").Value = 898
# Row 101
$ws.Range("H101").Value = 497.5
$ws.Range("I101").Value = 497.5
$ws.Range("K101").Value = 1492.5
$ws.Range("M101").Value = 129.5
# Row 105
$ws.Range("H105").Value = 38000
$ws.Range("J105").Value = 38000
$ws.Range("L105").Value = 38000
$ws.Range("N105").Value = -44988

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1507.6097
$ws.Range("I45").Value = 1495.3
$ws.Range("K45").Value = 1495.3
$ws.Range("M45").Value = -1118.3
# Row 97
$ws.Range("H97").Value = 1809.6923
$ws.Range("I97").Value = 1539.8334
$ws.Range("K97").Value = 1539.8334
$ws.Range("M97").Value = -1043.8334
# Row 102
$ws.Range("H102").Value = 1067.875
$ws.Range("I102").Value = 957.3333
$ws.Range("J102").Value = 1399.5
$ws.Range("K102").Value = 957.3333
$ws.Range("L102").Value = 1399.5
$ws.Range("M102").Value = 664.6667
$ws.Range("N102").Value = -4643.5
# Row 110
$ws.Range("H110").Value = 409.5
$ws.Range("I110").Value = 409.5
$ws.Range("K110").Value = 409.5
$ws.Range("M110").Value = 1635.5
# Row 131
$ws.Range("H131").Value = 7984.5
$ws.Range("I131").Value = 7000
$ws.Range("K131").Value = 21000
$ws.Range("M131").Value = -15960
# Row 138
$ws.Range("H138").Value = 2645.8235
$ws.Range("J138").Value = 3027.0715
$ws.Range("L138").Value = 9081.2145
$ws.Range("N138").Value = -19361.2145

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 2620.7778
$ws.Range("J107").Value = 3999.5
$ws.Range("L107").Value = 3999.5
$ws.Range("N107").Value = -7839.5
# Row 132
$ws.Range("H132").Value = 1578
$ws.Range("I132").Value = 937.3333
$ws.Range("K132").Value = 2811.9999
$ws.Range("M132").Value = -281.9998999999998

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 13
$ws.Range("I7").Value = 13
$ws.Range("K7").Value = 13
$ws.Range("M7").Value = 100
# Row 22
$ws.Range("H22").Value = 748.2353000000001
$ws.Range("I22").Value = 748.125
$ws.Range("K22").Value = 748.125
$ws.Range("M22").Value = -398.125
# Row 59
$ws.Range("H59").Value = 23000.334
$ws.Range("I59").Value = 11000.5
$ws.Range("J59").Value = 47000
$ws.Range("K59").Value = 11000.5
$ws.Range("L59").Value = 47000
$ws.Range("M59").Value = -9855.5
$ws.Range("N59").Value = -49290
# Row 60
$ws.Range("H60").Value = 27836.6
$ws.Range("I60").Value = 16061
$ws.Range("K60").Value = 16061
$ws.Range("M60").Value = -15550
# Row 69
$ws.Range("H69").Value = 10499.5
$ws.Range("I69").Value = 10499.5
$ws.Range("K69").Value = 10499.5
$ws.Range("M69").Value = -9750.5
# Row 72
$ws.Range("H72").Value = 10499.5
$ws.Range("I72").Value = 10499.5
$ws.Range("K72").Value = 31498.5
$ws.Range("M72").Value = -27754.5

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 2836.1943
$ws.Range("I4").Value = 2004.9524
$ws.Range("J4").Value = 3999.9333
$ws.Range("K4").Value = 6014.857199999999
$ws.Range("L4").Value = 11999.7999
$ws.Range("M4").Value = -5902.857199999999
$ws.Range("N4").Value = -12223.7999
# Row 34
$ws.Range("H34").Value = 718
$ws.Range("I34").Value = 392.33334
$ws.Range("J34").Value = 1695
$ws.Range("K34").Value = 1177.00002
$ws.Range("L34").Value = 5085
$ws.Range("M34").Value = -1093.00002
$ws.Range("N34").Value = -5253
# Row 63
$ws.Range("H63").Value = 3646.4443
# Row 66
$ws.Range("H66").Value = 3646.4443
# Row 75
$ws.Range("H75").Value = 2482.4
$ws.Range("J75").Value = 2482.4
$ws.Range("L75").Value = 7447.200000000001
$ws.Range("N75").Value = -9443.200000000001
# Row 78
$ws.Range("H78").Value = 2482.4
$ws.Range("J78").Value = 2482.4
$ws.Range("L78").Value = 22341.6
$ws.Range("N78").Value = -32325.6
# Row 93
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
# Row 107
$ws.Range("H107").Value = 499.5
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 499.5
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 1498.5
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -5338.5
# Row 109
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()
$ws.Range("N109").Value = -5080

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 1550
$ws.Range("I80").Value = 1500
$ws.Range("J80").Value = 1600
$ws.Range("K80").Value = 1500
$ws.Range("L80").Value = 1600
$ws.Range("M80").Value = -502
$ws.Range("N80").Value = -3596
# Row 83
$ws.Range("H83").Value = 1550
$ws.Range("I83").Value = 1500
$ws.Range("J83").Value = 1600
$ws.Range("K83").Value = 7500
$ws.Range("L83").Value = 8000
$ws.Range("M83").Value = -2508
$ws.Range("N83").Value = -17984
# Row 97
$ws.Range("H97").Value = 269
$ws.Range("I97").Value = 271.9091
$ws.Range("J97").Value = 258.33334
$ws.Range("K97").Value = 271.9091
$ws.Range("L97").Value = 258.33334
$ws.Range("M97").Value = 224.0909
$ws.Range("N97").Value = -1250.33334
# Row 106
$ws.Range("H106").Value = 25000
$ws.Range("J106").Value = 25000
$ws.Range("L106").Value = 25000
$ws.Range("N106").Value = -27524
# Row 113
$ws.Range("H113").Value = 1312.5
$ws.Range("I113").Value = 1312
$ws.Range("K113").Value = 1312
$ws.Range("M113").Value = 858

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 970.3333
$ws.Range("I7").Value = 970.3333
$ws.Range("K7").Value = 970.3333
$ws.Range("M7").Value = -858.3333
# Row 22
$ws.Range("H22").Value = 2133.8262
$ws.Range("I22").Value = 1052.4
$ws.Range("J22").Value = 4161.5
$ws.Range("K22").Value = 1052.4
$ws.Range("L22").Value = 4161.5
$ws.Range("M22").Value = -757.4000000000001
$ws.Range("N22").Value = -4751.5
# Row 27
$ws.Range("H27").Value = 2133.8262
$ws.Range("I27").Value = 1052.4
$ws.Range("J27").Value = 4161.5
$ws.Range("K27").Value = 1052.4
$ws.Range("L27").Value = 4161.5
$ws.Range("M27").Value = -945.4000000000001
$ws.Range("N27").Value = -4375.5
# Row 46
$ws.Range("H46").Value = 6747.625
$ws.Range("J46").Value = 6747.625
$ws.Range("L46").Value = 6747.625
$ws.Range("N46").Value = -7123.625
# Row 126
$ws.Range("H126").Value = 970.3333
$ws.Range("I126").Value = 970.3333
$ws.Range("K126").Value = 2910.9999
$ws.Range("M126").Value = -440.9998999999998

$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 28493.5
$ws.Range("I54").Value = 27989
$ws.Range("K54").Value = 27989
$ws.Range("M54").Value = -27469
# Row 62
$ws.Range("H62").Value = 4375
$ws.Range("I62").Value = 2500
$ws.Range("K62").Value = 2500
$ws.Range("M62").Value = -1876
# Row 65
$ws.Range("H65").Value = 4375
$ws.Range("I65").Value = 2500
$ws.Range("K65").Value = 12500
$ws.Range("M65").Value = -9380
# Row 104
$ws.Range("H104").Value = 18450
$ws.Range("J104").Value = 18450
$ws.Range("L104").Value = 18450
$ws.Range("N104").Value = -25438
# Row 113
$ws.Range("H113").Value = 3144.158
$ws.Range("I113").Value = 287.53845
$ws.Range("J113").Value = 9333.5
$ws.Range("K113").Value = 862.61535
$ws.Range("L113").Value = 28000.5
$ws.Range("M113").Value = 1307.38465
$ws.Range("N113").Value = -32340.5
